$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.590.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").Value = "'3.192.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'532.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.00%  "
$ws.Range("D6").Value = "'134.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.39%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'3.192.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.10%  "
$ws.Range("D9").Value = "'0.457"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.00%  "
$ws.Range("E10").Value = "  -6.68%  "
$ws.Range("E11").Value = "  -7.22%  "
$ws.Range("E12").Value = "  -5.56%  "
$ws.Range("D13").Value = "'3.741.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.28%  "
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "'25.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.98%  "
$ws.Range("D16").Value = "'3.197.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.96%  "
$ws.Range("D17").Value = "'58.656.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.87%  "
$ws.Range("E18").Value = "  -7.74%  "
$ws.Range("E19").Value = "  -6.93%  "
$ws.Range("D20").Value = "'13.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.91%  "
$ws.Range("E21").Value = "  -9.38%  "
$ws.Range("D22").Value = "'358.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.53%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'69.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.67%  "
$ws.Range("D25").Value = "'0.515"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.19%  "
$ws.Range("D26").Value = "'3.335.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.87%  "
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("D28").Value = "'0.0₃0950"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.75%  "
$ws.Range("D30").Value = "'7.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.06%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -8.35%  "
$ws.Range("D33").Value = "'6.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.13%  "
$ws.Range("D34").Value = "'21.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.58%  "
$ws.Range("E35").Value = "  -7.95%  "
$ws.Range("E36").Value = "  -4.98%  "
$ws.Range("D37").Value = "'4.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.79%  "
$ws.Range("D38").Value = "'6.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.10%  "
$ws.Range("D39").Value = "'1.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.76%  "
$ws.Range("D40").Value = "'25.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.92%  "
$ws.Range("D41").Value = "'0.0702"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.92%  "
$ws.Range("D42").Value = "'3.221.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("D43").Value = "'40.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("D44").Value = "'0.709"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.74%  "
$ws.Range("D45").Value = "'1.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("D46").Value = "'4.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.87%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.80%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "'2.283.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.35%  "
$ws.Range("D50").Value = "'6.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.66%  "
$ws.Range("D51").Value = "'20.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.02%  "
